$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cellUpdates = @(
    @{Cell='D2'; Val='67.368.28'},
    @{Cell='E2'; Val='  -4.68%  '},
    @{Cell='D3'; Val='3.267.57'},
    @{Cell='E3'; Val='  -7.18%  '},
    @{Cell='E4'; Val='  +0.02%  '},
    @{Cell='D5'; Val='592.22'},
    @{Cell='E5'; Val='  -4.69%  '},
    @{Cell='D6'; Val='151.76'},
    @{Cell='E6'; Val='  -11.63%  '},
    @{Cell='E7'; Val='  -0.06%  '},
    @{Cell='D8'; Val='3.258.44'},
    @{Cell='E8'; Val='  -7.39%  '},
    @{Cell='D9'; Val='0.543'},
    @{Cell='E9'; Val='  -10.69%  '},
    @{Cell='E10'; Val='  -13.79%  '},
    @{Cell='D11'; Val='6.60'},
    @{Cell='E11'; Val='  -7.93%  '},
    @{Cell='D12'; Val='0.514'},
    @{Cell='E12'; Val='  -12.04%  '},
    @{Cell='D13'; Val='38.59'},
    @{Cell='E13'; Val='  -16.44%  '},
    @{Cell='D14'; Val='0.0000245'},
    @{Cell='E14'; Val='  -11.04%  '},
    @{Cell='D15'; Val='3.788.57'},
    @{Cell='E15'; Val='  -7.16%  '},
    @{Cell='D16'; Val='67.342.38'},
    @{Cell='E16'; Val='  -4.89%  '},
    @{Cell='D17'; Val='3.265.95'},
    @{Cell='E17'; Val='  -6.90%  '},
    @{Cell='D18'; Val='7.27'},
    @{Cell='E18'; Val='  -13.76%  '},
    @{Cell='B19'; Val='BitcoinCash'},
    @{Cell='C19'; Val='https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'},
    @{Cell='D19'; Val='534.98'},
    @{Cell='E19'; Val='  -11.78%  '},
    @{Cell='B20'; Val='TRON'},
    @{Cell='C20'; Val='https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'},
    @{Cell='D20'; Val='0.114'},
    @{Cell='E20'; Val='  -6.10%  '},
    @{Cell='D21'; Val='15.16'},
    @{Cell='E21'; Val='  -14.25%  '},
    @{Cell='D22'; Val='0.764'},
    @{Cell='E22'; Val='  -13.17%  '},
    @{Cell='E23'; Val='  -12.72%  '},
    @{Cell='D24'; Val='85.90'},
    @{Cell='E24'; Val='  -11.69%  '},
    @{Cell='D25'; Val='13.65'},
    @{Cell='E25'; Val='  -12.27%  '},
    @{Cell='E26'; Val='  -0.01%  '},
    @{Cell='E27'; Val='  -12.31%  '},
    @{Cell='D28'; Val='8.11'},
    @{Cell='E28'; Val='  -10.13%  '},
    @{Cell='D29'; Val='29.44'},
    @{Cell='E29'; Val='  -12.24%  '},
    @{Cell='D30'; Val='2.15'},
    @{Cell='E30'; Val='  -15.92%  '},
    @{Cell='D31'; Val='2.68'},
    @{Cell='E31'; Val='  -10.57%  '},
    @{Cell='E32'; Val='  -11.64%  '},
    @{Cell='D33'; Val='545.20'},
    @{Cell='E33'; Val='  -11.70%  '},
    @{Cell='E34'; Val='  -17.97%  '},
    @{Cell='D35'; Val='5.77'},
    @{Cell='E35'; Val='  -15.12%  '},
    @{Cell='E36'; Val='  -0.03%  '},
    @{Cell='D37'; Val='0.0460'},
    @{Cell='E37'; Val='  -6.66%  '},
    @{Cell='D38'; Val='53.34'},
    @{Cell='E38'; Val='  -5.88%  '},
    @{Cell='D39'; Val='0.0864'},
    @{Cell='E39'; Val='  -13.00%  '},
    @{Cell='D40'; Val='9.13'},
    @{Cell='E40'; Val='  -15.87%  '},
    @{Cell='E41'; Val='  -10.12%  '},
    @{Cell='D42'; Val='2.76'},
    @{Cell='E42'; Val='  -18.59%  '},
    @{Cell='D43'; Val='2.946.70'},
    @{Cell='E43'; Val='  -11.66%  '},
    @{Cell='D44'; Val='0.0₃0595'},
    @{Cell='E44'; Val='  -17.66%  '},
    @{Cell='D45'; Val='0.264'},
    @{Cell='E45'; Val='  -14.66%  '},
    @{Cell='D46'; Val='26.91'},
    @{Cell='E46'; Val='  -15.38%  '},
    @{Cell='D47'; Val='2.17'},
    @{Cell='E47'; Val='  -12.92%  '},
    @{Cell='E48'; Val='  -0.02%  '},
    @{Cell='E49'; Val='  -19.53%  '},
    @{Cell='D50'; Val='126.86'},
    @{Cell='E50'; Val='  -5.22%  '},
    @{Cell='E51'; Val='  -12.31%  '}
)

foreach ($u in $cellUpdates) {
    $rng = $ws.Range($u.Cell)
    $rng.NumberFormat = "@"
    $rng.Value = $u.Val
    $rng.Style = "Normal"
}
